# Weekly update: insert 4 new price rows (latest week) at the top of the
# data block (row 41) for "Haba" / Mercado Mayorista Lo Valledor de
# Santiago, pushing all existing data rows down by 4. Sheet dimension
# grows from A1:R126 to A1:R130.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows right before the current row 41, shifting
# rows 41:126 down to 45:130.
$ws.Rows("41:44").Insert()

# Common (unchanged-across-all-four) field values for the newly inserted
# rows, matching the rest of the "Haba" block.
$mercadoId   = 6
$mercado     = "Mercado Mayorista Lo Valledor de Santiago"
$region      = "Metropolitana"
$codreg      = 13
$categoriaId = 100112026
$categoria   = "Haba"
$variedad    = "Sin especificar"
$unidad      = "`$/saco 25 kilos"
$kgUnidades  = 25
$clasif      = "Hortaliza"

# Per-row data: Fecha(serial), Calidad, Volumen, PrecioMin, PrecioMax,
# PrecioPromPonderado, Origen, Precio$/Kg
$rows = @(
    @{ Row = 41; Fecha = 44469; Calidad = "Primera"; Volumen = 1400; PMin = 6000;  PMax = 7000;  PProm = 6536; Origen = "Región Metropolitana"; PKg = 261 },
    @{ Row = 42; Fecha = 44469; Calidad = "Primera"; Volumen = 430;  PMin = 6000;  PMax = 7000;  PProm = 6465; Origen = "Región de Coquimbo";   PKg = 259 },
    @{ Row = 43; Fecha = 44469; Calidad = "Segunda"; Volumen = 400;  PMin = 5000;  PMax = 5000;  PProm = 5000; Origen = "Región Metropolitana"; PKg = 200 },
    @{ Row = 44; Fecha = 44469; Calidad = "Segunda"; Volumen = 130;  PMin = 5000;  PMax = 5000;  PProm = 5000; Origen = "Región de Coquimbo";   PKg = 200 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $ws.Cells.Item($rowNum, 1).Value  = $mercadoId
    $ws.Cells.Item($rowNum, 2).Value  = $mercado
    $ws.Cells.Item($rowNum, 3).Value  = $region
    $ws.Cells.Item($rowNum, 4).Value  = $r.Fecha
    $ws.Cells.Item($rowNum, 5).Value  = $codreg
    $ws.Cells.Item($rowNum, 6).Value  = $categoriaId
    $ws.Cells.Item($rowNum, 7).Value  = $categoria
    $ws.Cells.Item($rowNum, 8).Value  = $variedad
    $ws.Cells.Item($rowNum, 9).Value  = $r.Calidad
    $ws.Cells.Item($rowNum, 10).Value = $r.Volumen
    $ws.Cells.Item($rowNum, 11).Value = $r.PMin
    $ws.Cells.Item($rowNum, 12).Value = $r.PMax
    $ws.Cells.Item($rowNum, 13).Value = $r.PProm
    $ws.Cells.Item($rowNum, 14).Value = $unidad
    $ws.Cells.Item($rowNum, 15).Value = $r.Origen
    $ws.Cells.Item($rowNum, 16).Value = $r.PKg
    $ws.Cells.Item($rowNum, 17).Value = $kgUnidades
    $ws.Cells.Item($rowNum, 18).Value = $clasif
}
